$d = $word.ActiveDocument

# Step 1: Replace "type" paragraph (38) with "transactionType" + proofErr
$p38 = $d.Paragraphs(38)
$xml38 = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='3'/></w:numPr><w:rPr><w:sz w:val='24'/></w:rPr></w:pPr><w:proofErr w:type='spellStart'/><w:r><w:rPr><w:sz w:val='24'/></w:rPr><w:t>transactionType</w:t></w:r><w:proofErr w:type='spellEnd'/></w:p>"
$p38.Range.InsertXML($xml38)

# Step 2: Replace "userId" paragraph (39) with proofErr-wrapped userId, no bookmark
$p39 = $d.Paragraphs(39)
$xml39 = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='3'/></w:numPr><w:rPr><w:sz w:val='24'/></w:rPr></w:pPr><w:proofErr w:type='spellStart'/><w:r><w:rPr><w:sz w:val='24'/></w:rPr><w:t>userId</w:t></w:r><w:proofErr w:type='spellEnd'/></w:p>"
$p39.Range.InsertXML($xml39)

# Step 3: Append 3 new paragraphs after paragraph 39
$p39b = $d.Paragraphs(39)
$null = $p39b.Range.InsertParagraphAfter()
$p40 = $d.Paragraphs(40)
$null = $p40.Range.InsertParagraphAfter()
$p41 = $d.Paragraphs(41)
$null = $p41.Range.InsertParagraphAfter()

# Step 4: Fill in paragraph 40 (bookId) - no proofErr, WITH bookmark
$p40 = $d.Paragraphs(40)
$xml40 = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='3'/></w:numPr><w:rPr><w:sz w:val='24'/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val='24'/></w:rPr><w:t>bookId</w:t></w:r><w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/></w:p>"
$p40.Range.InsertXML($xml40)

# Step 5: Fill in paragraph 41 (pickupDate) - WITH proofErr
$p41 = $d.Paragraphs(41)
$xml41 = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='3'/></w:numPr><w:rPr><w:sz w:val='24'/></w:rPr></w:pPr><w:proofErr w:type='spellStart'/><w:r><w:rPr><w:sz w:val='24'/></w:rPr><w:t>pickupDate</w:t></w:r><w:proofErr w:type='spellEnd'/></w:p>"
$p41.Range.InsertXML($xml41)

# Step 6: Fill in paragraph 42 (rentalDate) - WITH proofErr
$p42 = $d.Paragraphs(42)
$xml42 = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='3'/></w:numPr><w:rPr><w:sz w:val='24'/></w:rPr></w:pPr><w:proofErr w:type='spellStart'/><w:r><w:rPr><w:sz w:val='24'/></w:rPr><w:t>rentalDate</w:t></w:r><w:proofErr w:type='spellEnd'/></w:p>"
$p42.Range.InsertXML($xml42)

Write-Host "total paragraphs: $($d.Paragraphs.Count)"
for ($i = 36; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs($i)
    Write-Host "$i : [$($pp.Range.Text)]"
}
